# Apply the cryptos-list price/volume update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '52.028.94'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.793.80'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.04%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '359.24'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '110.03'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.43%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.564'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E9').Value = '  -1.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.18'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.27%  '
$ws.Range('E12').Value = '  +1.28%  '
$ws.Range('E13').Value = '  -1.85%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.61'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.230.79'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.783.25'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.28%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.947'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +6.84%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '51.953.74'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.46%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.43'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.15%  '
$ws.Range('E20').Value = '  -1.72%  '
$ws.Range('E21').Value = '  -2.80%  '
$ws.Range('E22').Value = '  -0.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '273.75'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.32'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.84%  '
$ws.Range('E25').Value = '  -1.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.70'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.21%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.28'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.62%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.21'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.14%  '
$ws.Range('E30').Value = '  +4.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '51.65'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.83%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0464'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.57%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.48'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.51%  '
$ws.Range('E34').Value = '  -1.48%  '
$ws.Range('E35').Value = '  +2.94%  '
$ws.Range('E36').Value = '  -1.47%  '
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.24'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.29'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.51%  '
$ws.Range('E40').Value = '  -2.85%  '
$ws.Range('E41').Value = '  +2.20%  '
$ws.Range('E42').Value = '  -1.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '122.77'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.61%  '
$ws.Range('E44').Value = '  -2.34%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '22.19'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.32%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.088.01'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.27'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.09%  '
$ws.Range('E48').Value = '  -1.68%  '
$ws.Range('E49').Value = '  +1.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.936'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.32%  '
$ws.Range('E51').Value = '  +0.74%  '
